$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Clear existing hyperlinks (and their relationships) so we can rebuild them cleanly
$ws.Hyperlinks.Delete()

$ts = '2026-02-02 18:40:13'
$category = 'システム開発'
$deadline = '期限情報なし'

# Row 2
$ws.Cells.Item(2, 1).Value = $ts
$ws.Cells.Item(2, 2).Value = 'Gmail、スプレッドシート、Google Driveを連携した 業務効率化システム開発'
$ws.Cells.Item(2, 3).Value = $category
$ws.Cells.Item(2, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(2, 5).Value = $deadline
$ws.Cells.Item(2, 6).Value = 'https://www.lancers.jp/work/detail/5484339'
$ws.Cells.Item(2, 7).Value = 453
$ws.Cells.Item(2, 8).Value = '🔥AI,Ai ◆効率化,開発'
$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), 'https://www.lancers.jp/work/detail/5484339')
$ws.Cells.Item(2, 6).Style = "Hyperlink"

# Row 3
$ws.Cells.Item(3, 1).Value = $ts
$ws.Cells.Item(3, 2).Value = '【完全在宅】AI×Web開発エンジニア募集!業務自動化・AI機能開発'
$ws.Cells.Item(3, 3).Value = $category
$ws.Cells.Item(3, 4).Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Cells.Item(3, 5).Value = $deadline
$ws.Cells.Item(3, 6).Value = 'https://www.lancers.jp/work/detail/5483480'
$ws.Cells.Item(3, 7).Value = 435
$ws.Cells.Item(3, 8).Value = '🔥AI,Ai ◆開発,自動化'
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), 'https://www.lancers.jp/work/detail/5483480')
$ws.Cells.Item(3, 6).Style = "Hyperlink"

# Row 4
$ws.Cells.Item(4, 1).Value = $ts
$ws.Cells.Item(4, 2).Value = '製造業向け設備要件定義書の自動生成AIシステムの開発・DB設計支援エンジニア(AI/バックエンド)'
$ws.Cells.Item(4, 3).Value = $category
$ws.Cells.Item(4, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(4, 5).Value = $deadline
$ws.Cells.Item(4, 6).Value = 'https://www.lancers.jp/work/detail/5473648'
$ws.Cells.Item(4, 7).Value = 390
$ws.Cells.Item(4, 8).Value = '🔥AI,Ai ◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), 'https://www.lancers.jp/work/detail/5473648')
$ws.Cells.Item(4, 6).Style = "Hyperlink"

# Row 5
$ws.Cells.Item(5, 1).Value = $ts
$ws.Cells.Item(5, 2).Value = '※急募:Next.jsによる業務アプリの開発(+Flutter)'
$ws.Cells.Item(5, 3).Value = $category
$ws.Cells.Item(5, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(5, 5).Value = $deadline
$ws.Cells.Item(5, 6).Value = 'https://www.lancers.jp/work/detail/5483966'
$ws.Cells.Item(5, 7).Value = 225
$ws.Cells.Item(5, 8).Value = '🔥Next.js ◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), 'https://www.lancers.jp/work/detail/5483966')
$ws.Cells.Item(5, 6).Style = "Hyperlink"

# Row 6
$ws.Cells.Item(6, 1).Value = $ts
$ws.Cells.Item(6, 2).Value = '※急募:Flutterによる業務アプリの開発(+next.js)'
$ws.Cells.Item(6, 3).Value = $category
$ws.Cells.Item(6, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(6, 5).Value = $deadline
$ws.Cells.Item(6, 6).Value = 'https://www.lancers.jp/work/detail/5483967'
$ws.Cells.Item(6, 7).Value = 218
$ws.Cells.Item(6, 8).Value = '🔥Next.js ◆開発 ◇アプリ'
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), 'https://www.lancers.jp/work/detail/5483967')
$ws.Cells.Item(6, 6).Style = "Hyperlink"

# Row 7
$ws.Cells.Item(7, 1).Value = $ts
$ws.Cells.Item(7, 2).Value = '【急募】楽天RPP広告自動化ツールの開発依頼'
$ws.Cells.Item(7, 3).Value = $category
$ws.Cells.Item(7, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(7, 5).Value = $deadline
$ws.Cells.Item(7, 6).Value = 'https://www.lancers.jp/work/detail/5484001'
$ws.Cells.Item(7, 7).Value = 213
$ws.Cells.Item(7, 8).Value = '◆ツール,開発'
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), 'https://www.lancers.jp/work/detail/5484001')
$ws.Cells.Item(7, 6).Style = "Hyperlink"

# Row 8
$ws.Cells.Item(8, 1).Value = $ts
$ws.Cells.Item(8, 2).Value = 'X(旧twitter)のロック解除自動化システム構築'
$ws.Cells.Item(8, 3).Value = $category
$ws.Cells.Item(8, 4).Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Cells.Item(8, 5).Value = $deadline
$ws.Cells.Item(8, 6).Value = 'https://www.lancers.jp/work/detail/5484231'
$ws.Cells.Item(8, 7).Value = 103
$ws.Cells.Item(8, 8).Value = '◆自動化'
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), 'https://www.lancers.jp/work/detail/5484231')
$ws.Cells.Item(8, 6).Style = "Hyperlink"

# Row 9
$ws.Cells.Item(9, 1).Value = $ts
$ws.Cells.Item(9, 2).Value = '(仕様削減)【受注メールを元にしたスクレピング&抽出情報管理ツール】'
$ws.Cells.Item(9, 3).Value = $category
$ws.Cells.Item(9, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(9, 5).Value = $deadline
$ws.Cells.Item(9, 6).Value = 'https://www.lancers.jp/work/detail/5484248'
$ws.Cells.Item(9, 7).Value = 98
$ws.Cells.Item(9, 8).Value = '◆ツール ◇管理'
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), 'https://www.lancers.jp/work/detail/5484248')
$ws.Cells.Item(9, 6).Style = "Hyperlink"

# Row 10
$ws.Cells.Item(10, 1).Value = $ts
$ws.Cells.Item(10, 2).Value = '製造業DXプロダクト開発のプロダクトマネージャー募集'
$ws.Cells.Item(10, 3).Value = $category
$ws.Cells.Item(10, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(10, 5).Value = $deadline
$ws.Cells.Item(10, 6).Value = 'https://www.lancers.jp/work/detail/5468432'
$ws.Cells.Item(10, 7).Value = 75
$ws.Cells.Item(10, 8).Value = '◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), 'https://www.lancers.jp/work/detail/5468432')
$ws.Cells.Item(10, 6).Style = "Hyperlink"

# Row 11
$ws.Cells.Item(11, 1).Value = $ts
$ws.Cells.Item(11, 2).Value = '【急募】新しいWebサービスの開発パートナーを探しています!'
$ws.Cells.Item(11, 3).Value = $category
$ws.Cells.Item(11, 4).Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Cells.Item(11, 5).Value = $deadline
$ws.Cells.Item(11, 6).Value = 'https://www.lancers.jp/work/detail/5483482'
$ws.Cells.Item(11, 7).Value = 75
$ws.Cells.Item(11, 8).Value = '◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), 'https://www.lancers.jp/work/detail/5483482')
$ws.Cells.Item(11, 6).Style = "Hyperlink"

# Row 12
$ws.Cells.Item(12, 1).Value = $ts
$ws.Cells.Item(12, 2).Value = '四柱推命の命式自動計算プログラム(Web/Excel)の開発依頼'
$ws.Cells.Item(12, 3).Value = $category
$ws.Cells.Item(12, 4).Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Cells.Item(12, 5).Value = $deadline
$ws.Cells.Item(12, 6).Value = 'https://www.lancers.jp/work/detail/5484177'
$ws.Cells.Item(12, 7).Value = 68
$ws.Cells.Item(12, 8).Value = '◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), 'https://www.lancers.jp/work/detail/5484177')
$ws.Cells.Item(12, 6).Style = "Hyperlink"

# Row 13
$ws.Cells.Item(13, 1).Value = $ts
$ws.Cells.Item(13, 2).Value = '美容皮膚科向け LINE連携型BtoB SaaS(MVP) の開発案件'
$ws.Cells.Item(13, 3).Value = $category
$ws.Cells.Item(13, 4).Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Cells.Item(13, 5).Value = $deadline
$ws.Cells.Item(13, 6).Value = 'https://www.lancers.jp/work/detail/5483503'
$ws.Cells.Item(13, 7).Value = 68
$ws.Cells.Item(13, 8).Value = '◆開発'
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), 'https://www.lancers.jp/work/detail/5483503')
$ws.Cells.Item(13, 6).Style = "Hyperlink"

# Row 14
$ws.Cells.Item(14, 1).Value = $ts
$ws.Cells.Item(14, 2).Value = '【急募】Notion×Slackでのオンライン講座運営システム構築'
$ws.Cells.Item(14, 3).Value = $category
$ws.Cells.Item(14, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(14, 5).Value = $deadline
$ws.Cells.Item(14, 6).Value = 'https://www.lancers.jp/work/detail/5483854'
$ws.Cells.Item(14, 7).Value = 28
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), 'https://www.lancers.jp/work/detail/5483854')
$ws.Cells.Item(14, 6).Style = "Hyperlink"

# Row 15
$ws.Cells.Item(15, 1).Value = $ts
$ws.Cells.Item(15, 2).Value = '【市場調査】海外向けデジタルサービスの価値評価依頼'
$ws.Cells.Item(15, 3).Value = $category
$ws.Cells.Item(15, 4).Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Cells.Item(15, 5).Value = $deadline
$ws.Cells.Item(15, 6).Value = 'https://www.lancers.jp/work/detail/5483504'
$ws.Cells.Item(15, 7).Value = 13
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), 'https://www.lancers.jp/work/detail/5483504')
$ws.Cells.Item(15, 6).Style = "Hyperlink"

# Column B width: 46 -> 51 characters (COM ColumnWidth excludes ~0.833 cell padding)
$ws.Columns.Item(2).ColumnWidth = 50.166667

$ws.Range("A1").Select()
